# Updated symbol list on Thu Feb 16 11:27:09 UTC 2023 with GitHub Actions
# Refresh Price (column D) and Volume(1h) (column E) figures for the
# crypto rows that moved. Values are stored as text (leading apostrophe
# forces Excel to keep them as literal text, matching the sheet's
# existing text-formatted Price/Volume columns) so things like trailing
# zeros ("48.30") and percent signs are preserved exactly.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'321.06"
$ws.Range("E2").Value = "'7.22%"
$ws.Range("D3").Value = "'48.30"
$ws.Range("E3").Value = "'14.41%"
$ws.Range("D4").Value = "'5.255"
$ws.Range("E4").Value = "'4.80%"
$ws.Range("E5").Value = "'7.01%"
$ws.Range("D6").Value = "'4.596"
$ws.Range("E6").Value = "'4.89%"
$ws.Range("D7").Value = "'1.646"
$ws.Range("E7").Value = "'2.86%"
$ws.Range("D8").Value = "'1.209"
$ws.Range("E8").Value = "'28.49%"
$ws.Range("E9").Value = "'8.70%"
$ws.Range("D10").Value = "'0.1942"
$ws.Range("E10").Value = "'5.73%"
$ws.Range("D11").Value = "'0.09454"
$ws.Range("E11").Value = "'3.74%"
$ws.Range("D12").Value = "'0.04607"
$ws.Range("E12").Value = "'9.73%"
$ws.Range("D13").Value = "'0.1049"
$ws.Range("E13").Value = "'-0.03%"
$ws.Range("D14").Value = "'0.001334"
$ws.Range("E14").Value = "'3.70%"
$ws.Range("D15").Value = "'0.04168"
$ws.Range("E15").Value = "'1.55%"
$ws.Range("D16").Value = "'0.005874"
$ws.Range("E16").Value = "'1.59%"
$ws.Range("D17").Value = "'3.337"
$ws.Range("E17").Value = "'-0.20%"
$ws.Range("D18").Value = "'2.427"
$ws.Range("E18").Value = "'1.78%"
$ws.Range("D19").Value = "'0.3387"
$ws.Range("E19").Value = "'1.57%"
$ws.Range("D20").Value = "'8.073"
$ws.Range("E20").Value = "'-3.73%"
$ws.Range("D21").Value = "'0.1381"
$ws.Range("E21").Value = "'-1.97%"
$ws.Range("D22").Value = "'0.3125"
$ws.Range("E22").Value = "'-5.23%"
$ws.Range("D23").Value = "'0.001305"
$ws.Range("E23").Value = "'3.17%"
$ws.Range("E24").Value = "'9.04%"
$ws.Range("D25").Value = "'0.0001351"
$ws.Range("E25").Value = "'6.47%"
$ws.Range("D26").Value = "'0.0003541"
$ws.Range("E26").Value = "'-4.91%"
$ws.Range("D38").Value = "'0.02714"
$ws.Range("E38").Value = "'12.54%"
$ws.Range("D39").Value = "'0.05745"
$ws.Range("E39").Value = "'9.46%"
$ws.Range("D40").Value = "'0.006303"
$ws.Range("E40").Value = "'-5.59%"
$ws.Range("D41").Value = "'0.007869"
$ws.Range("E41").Value = "'2.23%"
$ws.Range("D42").Value = "'0.1440"
$ws.Range("E42").Value = "'8.37%"
$ws.Range("D43").Value = "'0.007707"
$ws.Range("E43").Value = "'4.36%"
$ws.Range("E44").Value = "'3.85%"
$ws.Range("E45").Value = "'6.29%"
$ws.Range("D46").Value = "'0.00006898"
$ws.Range("E46").Value = "'10.49%"
$ws.Range("D47").Value = "'0.00000000750"
$ws.Range("E47").Value = "'0.11%"
$ws.Range("D48").Value = "'0.06622"
$ws.Range("E48").Value = "'46.37%"
$ws.Range("D49").Value = "'0.004002"
$ws.Range("E49").Value = "'-4.73%"
$ws.Range("D50").Value = "'0.00002101"
$ws.Range("E50").Value = "'0.11%"
$ws.Range("D51").Value = "'0.0002001"
$ws.Range("E51").Value = "'0.11%"
